$wb = $excel.ActiveWorkbook

# --- Rechnung2: add an (optional) invoice-date row above the items table ---
$ws2 = $wb.Worksheets.Item("Rechnung2")

# Make room: push the header row (and everything below it) down by one row.
$ws2.Rows("9:9").Insert()

# Row 8 was already blank and stays in place; fill it with the new fields.
$ws2.Range("A8").Value = "Rechnungsdatum:"

# Give the date cell the same date format/style as the other date cells
# in the table (column B) before writing the date value into it.
$ws2.Range("B11").Copy()
$ws2.Range("B8").PasteSpecial(-4122)
$ws2.Range("B8").Value = 45321

$ws2.Range("C8").Value = "(optional)"

$ws2.Activate()
$ws2.Range("C8").Select()

# --- Tabelle1: document the new "Rechnungsdatum:" field ---
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Row 5 is already blank; insert a new blank row below it (pushing the
# header row and everything after it down by one) ...
$ws1.Rows("6:6").Insert()

# ... then fill row 5 with the explanation text.
$ws1.Range("A5").Value = "Rechnungsdatum: "
$ws1.Range("B5").Value = 'Das Datum für die Rechnung; Es wird nach "Rechnungsdatum:" gesucht.'

# --- Final selection state: user ends up on Tabelle1 at B5 ---
$ws1.Activate()
$ws1.Range("B5").Select()
